{"js": "// Load all body paragraphs once; indices are stable across the edits below\n// because we only change run/text content (via Replace) and insert one new\n// paragraph at the very end of the sequence of edits.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// --- Change 1: TestCase02_T1.sql -> TestCase03_T1.sql -------------------\nparagraphs.items[3].insertText(\n  \"Giao t\u00e1c \u201c\u0110\u0103ng k\u00fd \u0111\u1ed3 \u00e1n\u201d tham chi\u1ebfu \u0111\u1ebfn TestCase03_T1.sql.\",\n  \"Replace\"\n);\n\n// --- Change 2: TestCase02_T2.sql -> TestCase03_T2.sql -------------------\nparagraphs.items[4].insertText(\n  \"Giao t\u00e1c \u201cC\u1eadp nh\u1eadt th\u1eddi gian n\u1ed9p \u0111\u1ed3 \u00e1n\u201d tham chi\u1ebfu \u0111\u1ebfn TestCase03_T2.sql\",\n  \"Replace\"\n);\n\n// --- Change 3: T1 reads & saves the submission time, add assumed dates --\nparagraphs.items[7].insertText(\n  \"Giao t\u00e1c T1 \u0111\u1ecdc th\u1eddi gian n\u1ed9p c\u1ee7a \u0111\u1ed3 \u00e1n v\u00e0 l\u01b0u l\u1ea1i. Gi\u1ea3 s\u1eed th\u1eddi gian n\u1ed9p l\u00fac n\u00e0y l\u00e0 30/07/2012, ng\u00e0y k\u1ebft th\u00fac m\u00f4n h\u1ecdc \u1ee9ng v\u1edbi \u0111\u1ed3 \u00e1n tr\u00ean l\u00e0 01/08/2012.\",\n  \"Replace\"\n);\n\n// --- Change 4: T1 checks validity \"\u0111\u1ec3 th\u00eam d\u1eef li\u1ec7u ...\" -----------------\nparagraphs.items[8].insertText(\n  \"Giao t\u00e1c T1 ki\u1ec3m tra th\u1eddi gian n\u1ed9p \u0111\u1ed3 \u00e1n xem c\u00f3 h\u1ee3p l\u1ec7 kh\u00f4ng \u0111\u1ec3 th\u00eam d\u1eef li\u1ec7u v\u00e0o b\u1ea3ng DE_SINHVIEN.\",\n  \"Replace\"\n);\n\n// --- Change 5: T2 changes the date to 29/07/2012 (was 29/05/2012) -------\nparagraphs.items[9].insertText(\n  \"Ngay l\u00fac n\u00e0y giao t\u00e1c T2 thay \u0111\u1ed5i th\u1eddi gian c\u1ee7a \u0111\u1ed3 \u00e1n (m\u00e3 s\u1ed1 2) th\u00e0nh 29/07/2012.\",\n  \"Replace\"\n);\n\nawait context.sync();\n\n// --- Change 6: insert a brand-new paragraph before \"Giao t\u00e1c T1 \u0111\u1ecdc l\u1ea1i...\"\n// Re-fetch paragraphs since a new one will be added; items[10] is still the\n// \"\u0111\u1ecdc l\u1ea1i\" paragraph (edits above only replaced text in place, they did not\n// add/remove paragraphs).\nconst paragraphs2 = context.document.body.paragraphs;\nparagraphs2.load(\"items/text\");\nawait context.sync();\n\nconst target = paragraphs2.items[10];\ntarget.insertParagraph(\n  \"Sau \u0111\u00f3 T1 th\u1ea5y d\u1eef li\u1ec7u h\u1ee3p l\u1ec7 v\u00e0 xu\u1ea5t ra m\u00e0n h\u00ecnh th\u1eddi h\u1ea1n n\u1ed9p \u0111\u00e3 l\u01b0u tr\u01b0\u1edbc \u0111\u00f3 (b\u1eb1ng l\u1ec7nh PRINT) l\u00e0 30/07/2012.\",\n  \"Before\"\n);\nawait context.sync();\n\n// --- Change 7: update the \"\u0111\u00e3 tr\u1edf th\u00e0nh 29/05/2012\" paragraph to 29/07/2012\n// After the insertion above, the \"\u0111\u1ecdc l\u1ea1i\" paragraph shifted down by one.\nconst paragraphs3 = context.document.body.paragraphs;\nparagraphs3.load(\"items/text\");\nawait context.sync();\n\nparagraphs3.items[11].insertText(\n  \"Giao t\u00e1c T1 \u0111\u1ecdc l\u1ea1i th\u1eddi gian n\u1ed9p b\u1eb1ng l\u1ec7nh SELECT. Th\u1eddi gian n\u1ed9p l\u00fac n\u00e0y \u0111\u00e3 tr\u1edf th\u00e0nh  29/07/2012.\",\n  \"Replace\"\n);\n\n// --- Change 8: conclusion sentence rewritten -----------------------------\nparagraphs3.items[12].insertText(\n  \"K\u1ebft lu\u1eadn: th\u1eddi gian n\u1ed9p do giao t\u00e1c T1 \u0111\u1ecdc \u1edf hai l\u1ea7n (tr\u01b0\u1edbc v\u00e0 sau khi T2 ch\u1ea1y) c\u00f3 gi\u00e1 tr\u1ecb kh\u00e1c nhau.\",\n  \"Replace\"\n);\n\n// --- Change 9: fix sentence expanded with two more sentences ------------\nparagraphs3.items[13].insertText(\n  \"C\u00e1ch kh\u1eafc ph\u1ee5c: s\u1eed d\u1ee5ng m\u1ee9c c\u00f4 l\u1eadp \u201crepeatable read\u201d tr\u00ean giao t\u00e1c T1. Khi \u0111\u00f3 T1 tr\u01b0\u1edbc khi \u0111\u1ecdc s\u1ebd xin kh\u00f3a v\u00e0 kh\u00f3a \u0111\u01b0\u1ee3c gi\u1eef \u0111\u1ebfn h\u1ebft giao t\u00e1c n\u00e0y. Giao t\u00e1c T2 s\u1ebd ch\u1edd T1 tr\u1ea3 kh\u00f3a sau \u0111\u00f3 m\u1edbi th\u1ef1c hi\u1ec7n vi\u1ec7c c\u1eadp nh\u1eadt th\u1eddi gian n\u1ed9p. Do v\u1eady khi T1 xu\u1ea5t k\u1ebft qu\u1ea3 th\u1eddi gian n\u1ed9p \u1edf hai l\u1ea7n s\u1ebd ra gi\u1ed1ng nhau l\u00e0 30/07/2012. Sau khi T1 k\u1ebft th\u00fac, T2 m\u1edbi \u0111\u01b0\u1ee3c ph\u00e9p c\u1eadp nh\u1eadt th\u1eddi gian n\u1ed9p th\u00e0nh 29/07/2012.\",\n  \"Replace\"\n);\n\nawait context.sync();\n", "ps1": "# Replace a whole paragraph's visible text while leaving its paragraph\n# mark (and therefore its paragraph formatting / numbering) untouched.\nfunction Set-ParaText {\n    param($para, [string]$text)\n    $r = $para.Range\n    $r.End = $r.End - 1\n    $r.Text = $text\n}\n\n$d = $word.ActiveDocument\n\n# --- Change 1: TestCase02_T1.sql -> TestCase03_T1.sql -------------------\nSet-ParaText $d.Paragraphs.Item(4) \"Giao t\u00e1c \u201c\u0110\u0103ng k\u00fd \u0111\u1ed3 \u00e1n\u201d tham chi\u1ebfu \u0111\u1ebfn TestCase03_T1.sql.\"\n\n# --- Change 2: TestCase02_T2.sql -> TestCase03_T2.sql -------------------\nSet-ParaText $d.Paragraphs.Item(5) \"Giao t\u00e1c \u201cC\u1eadp nh\u1eadt th\u1eddi gian n\u1ed9p \u0111\u1ed3 \u00e1n\u201d tham chi\u1ebfu \u0111\u1ebfn TestCase03_T2.sql\"\n\n# --- Change 3: T1 reads & saves the submission time, add assumed dates --\nSet-ParaText $d.Paragraphs.Item(8) \"Giao t\u00e1c T1 \u0111\u1ecdc th\u1eddi gian n\u1ed9p c\u1ee7a \u0111\u1ed3 \u00e1n v\u00e0 l\u01b0u l\u1ea1i. Gi\u1ea3 s\u1eed th\u1eddi gian n\u1ed9p l\u00fac n\u00e0y l\u00e0 30/07/2012, ng\u00e0y k\u1ebft th\u00fac m\u00f4n h\u1ecdc \u1ee9ng v\u1edbi \u0111\u1ed3 \u00e1n tr\u00ean l\u00e0 01/08/2012.\"\n\n# --- Change 4: T1 checks validity \"\u0111\u1ec3 th\u00eam d\u1eef li\u1ec7u ...\" -----------------\nSet-ParaText $d.Paragraphs.Item(9) \"Giao t\u00e1c T1 ki\u1ec3m tra th\u1eddi gian n\u1ed9p \u0111\u1ed3 \u00e1n xem c\u00f3 h\u1ee3p l\u1ec7 kh\u00f4ng \u0111\u1ec3 th\u00eam d\u1eef li\u1ec7u v\u00e0o b\u1ea3ng DE_SINHVIEN.\"\n\n# --- Change 5: T2 changes the date to 29/07/2012 (was 29/05/2012) -------\nSet-ParaText $d.Paragraphs.Item(10) \"Ngay l\u00fac n\u00e0y giao t\u00e1c T2 thay \u0111\u1ed5i th\u1eddi gian c\u1ee7a \u0111\u1ed3 \u00e1n (m\u00e3 s\u1ed1 2) th\u00e0nh 29/07/2012.\"\n\n# --- Change 6: insert a brand-new paragraph before \"Giao t\u00e1c T1 \u0111\u1ecdc l\u1ea1i...\"\n# (currently paragraph 11: \"Giao t\u00e1c T1 \u0111\u1ecdc l\u1ea1i th\u1eddi gian n\u1ed9p ...\")\n$target = $d.Paragraphs.Item(11)\n$target.Range.InsertParagraphBefore()\nSet-ParaText $d.Paragraphs.Item(11) \"Sau \u0111\u00f3 T1 th\u1ea5y d\u1eef li\u1ec7u h\u1ee3p l\u1ec7 v\u00e0 xu\u1ea5t ra m\u00e0n h\u00ecnh th\u1eddi h\u1ea1n n\u1ed9p \u0111\u00e3 l\u01b0u tr\u01b0\u1edbc \u0111\u00f3 (b\u1eb1ng l\u1ec7nh PRINT) l\u00e0 30/07/2012.\"\n\n# --- Change 7: update \"\u0111\u00e3 tr\u1edf th\u00e0nh 29/05/2012\" paragraph to 29/07/2012 -\n# (shifted down to 12 after the insertion above)\nSet-ParaText $d.Paragraphs.Item(12) \"Giao t\u00e1c T1 \u0111\u1ecdc l\u1ea1i th\u1eddi gian n\u1ed9p b\u1eb1ng l\u1ec7nh SELECT. Th\u1eddi gian n\u1ed9p l\u00fac n\u00e0y \u0111\u00e3 tr\u1edf th\u00e0nh  29/07/2012.\"\n\n# --- Change 8: conclusion sentence rewritten -----------------------------\nSet-ParaText $d.Paragraphs.Item(13) \"K\u1ebft lu\u1eadn: th\u1eddi gian n\u1ed9p do giao t\u00e1c T1 \u0111\u1ecdc \u1edf hai l\u1ea7n (tr\u01b0\u1edbc v\u00e0 sau khi T2 ch\u1ea1y) c\u00f3 gi\u00e1 tr\u1ecb kh\u00e1c nhau.\"\n\n# --- Change 9: fix sentence expanded with two more sentences ------------\nSet-ParaText $d.Paragraphs.Item(14) \"C\u00e1ch kh\u1eafc ph\u1ee5c: s\u1eed d\u1ee5ng m\u1ee9c c\u00f4 l\u1eadp \u201crepeatable read\u201d tr\u00ean giao t\u00e1c T1. Khi \u0111\u00f3 T1 tr\u01b0\u1edbc khi \u0111\u1ecdc s\u1ebd xin kh\u00f3a v\u00e0 kh\u00f3a \u0111\u01b0\u1ee3c gi\u1eef \u0111\u1ebfn h\u1ebft giao t\u00e1c n\u00e0y. Giao t\u00e1c T2 s\u1ebd ch\u1edd T1 tr\u1ea3 kh\u00f3a sau \u0111\u00f3 m\u1edbi th\u1ef1c hi\u1ec7n vi\u1ec7c c\u1eadp nh\u1eadt th\u1eddi gian n\u1ed9p. Do v\u1eady khi T1 xu\u1ea5t k\u1ebft qu\u1ea3 th\u1eddi gian n\u1ed9p \u1edf hai l\u1ea7n s\u1ebd ra gi\u1ed1ng nhau l\u00e0 30/07/2012. Sau khi T1 k\u1ebft th\u00fac, T2 m\u1edbi \u0111\u01b0\u1ee3c ph\u00e9p c\u1eadp nh\u1eadt th\u1eddi gian n\u1ed9p th\u00e0nh 29/07/2012.\"\n"}
